$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update four odds values in row 4 (Hobro vs Hillerod)
$ws.Range("O4").Value = 1.18
$ws.Range("P4").Value = 4.5
$ws.Range("Q4").Value = 1.62
$ws.Range("R4").Value = 2.25

# 2. Delete the entire row 6 (Gloria Buzau vs Petrolul) -
#    this shifts rows 7 and 8 up by one, matching the diff.
$ws.Rows.Item(6).Delete()
